$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three account ("Conta") numbers whose rows must be removed entirely
# from the "Export" sheet (column A holds the account number as text).
$accountsToRemove = @("004946997", "004550415", "004431689")

# Locate the current row number for each account via Find (robust to the
# exact row position), collect them, then delete from the bottom row
# upwards so earlier row numbers remain valid while we delete.
$rowsToDelete = @()
foreach ($acct in $accountsToRemove) {
    $found = $ws.Columns.Item(1).Find($acct)
    $rowsToDelete += $found.Row
}

$rowsToDelete = $rowsToDelete | Sort-Object -Descending

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
